$d = $word.ActiveDocument

$replacements = @(
    @{old="725×8="; new="309×8="},
    @{old="357×4="; new="854×6="},
    @{old="443×4="; new="171×6="},
    @{old="899×9="; new="990×4="},
    @{old="597×3="; new="184×9="},
    @{old="569×4="; new="959×3="},
    @{old="239×5="; new="561×6="},
    @{old="274×6="; new="629×2="},
    @{old="498×9="; new="233×3="},
    @{old="960×7="; new="163×4="},
    @{old="669×2="; new="914×5="},
    @{old="205×4="; new="269×3="},
    @{old="372×8="; new="795×8="},
    @{old="443×7="; new="878×3="},
    @{old="163×9="; new="153×4="},
    @{old="329×8="; new="765×9="},
    @{old="580×3="; new="983×2="},
    @{old="406×2="; new="305×7="},
    @{old="535×6="; new="381×3="},
    @{old="666×5="; new="410×6="},
    @{old="289×3="; new="348×7="},
    @{old="989×6="; new="645×6="},
    @{old="868×5="; new="175×3="},
    @{old="490×5="; new="277×4="},
    @{old="495×4="; new="617×4="}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
